# Insert a new weekly data row at row 7 (pushing the existing rows 7-52
# down to 8-53) and populate it with the new "Fruta / hortaliza, semanal"
# observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 7..52 down to 8..53, creating a blank row 7 that inherits the
# formatting (e.g. the date style on column D) of the row it displaces.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new record.
$ws.Cells.Item(7, 1).Value = 8
$ws.Cells.Item(7, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(7, 3).Value = 'Coquimbo'
$ws.Cells.Item(7, 4).Value = 45282
$ws.Cells.Item(7, 5).Value = 4
$ws.Cells.Item(7, 6).Value = 'Fruta'
$ws.Cells.Item(7, 7).Value = 100101
$ws.Cells.Item(7, 8).Value = 'Berries'
$ws.Cells.Item(7, 9).Value = 100101001
$ws.Cells.Item(7, 10).Value = 'Arándano (blue)'
$ws.Cells.Item(7, 11).Value = 'Sin especificar'
$ws.Cells.Item(7, 12).Value = 'Primera'
$ws.Cells.Item(7, 13).Value = 400
$ws.Cells.Item(7, 14).Value = 7500
$ws.Cells.Item(7, 15).Value = 8000
$ws.Cells.Item(7, 16).Value = 7750
$ws.Cells.Item(7, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(7, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(7, 19).Value = 3875
$ws.Cells.Item(7, 20).Value = 2
